# Insert two new data rows (133 and 134) into the "Durazno" sheet.
# Inserting a 2-row range shifts all existing rows (old 133..205) down by
# two positions (to 135..207), which reproduces every cascading change
# seen across the rest of the diff (K/L/M/N/O/P/Q/R/S/T values "moving"
# down two rows). We only need to populate the two freshly-inserted rows
# with their brand-new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 133 downwards by inserting a fresh 2-row block before them.
$ws.Range("A133:T134").Insert()

# New row 133: Durazno / Kurakata / Primera
$ws.Cells.Item(133, 1).Value2  = 7
$ws.Cells.Item(133, 2).Value2  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(133, 3).Value2  = "Ñuble"
$ws.Cells.Item(133, 4).Value2  = 44572
$ws.Cells.Item(133, 5).Value2  = 16
$ws.Cells.Item(133, 6).Value2  = "Fruta"
$ws.Cells.Item(133, 7).Value2  = 100103
$ws.Cells.Item(133, 8).Value2  = "Frutos de hueso (carozo)"
$ws.Cells.Item(133, 9).Value2  = 100103004
$ws.Cells.Item(133, 10).Value2 = "Durazno"
$ws.Cells.Item(133, 11).Value2 = "Kurakata"
$ws.Cells.Item(133, 12).Value2 = "Primera"
$ws.Cells.Item(133, 13).Value2 = 200
$ws.Cells.Item(133, 14).Value2 = 13000
$ws.Cells.Item(133, 15).Value2 = 14000
$ws.Cells.Item(133, 16).Value2 = 13500
$ws.Cells.Item(133, 17).Value2 = "$/caja 16 kilos empedrada"
$ws.Cells.Item(133, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(133, 19).Value2 = 844
$ws.Cells.Item(133, 20).Value2 = 16

# New row 134: Durazno / Kurakata / Segunda
$ws.Cells.Item(134, 1).Value2  = 7
$ws.Cells.Item(134, 2).Value2  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(134, 3).Value2  = "Ñuble"
$ws.Cells.Item(134, 4).Value2  = 44572
$ws.Cells.Item(134, 5).Value2  = 16
$ws.Cells.Item(134, 6).Value2  = "Fruta"
$ws.Cells.Item(134, 7).Value2  = 100103
$ws.Cells.Item(134, 8).Value2  = "Frutos de hueso (carozo)"
$ws.Cells.Item(134, 9).Value2  = 100103004
$ws.Cells.Item(134, 10).Value2 = "Durazno"
$ws.Cells.Item(134, 11).Value2 = "Kurakata"
$ws.Cells.Item(134, 12).Value2 = "Segunda"
$ws.Cells.Item(134, 13).Value2 = 120
$ws.Cells.Item(134, 14).Value2 = 11000
$ws.Cells.Item(134, 15).Value2 = 12000
$ws.Cells.Item(134, 16).Value2 = 11500
$ws.Cells.Item(134, 17).Value2 = "$/caja 16 kilos empedrada"
$ws.Cells.Item(134, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(134, 19).Value2 = 719
$ws.Cells.Item(134, 20).Value2 = 16
